$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sources")
$ws.Activate()

# --- Fix the "ESS 2023" entry (B11): year typo 2024 -> 2023, keep the
# italic formatting on the "Coming Summer ..." run intact.
$essCell = $ws.Range("B11")
$essCell.Characters(26, 4).Text = "2023"
$essCell.Characters(12, 19).Font.Italic = $true

# --- Fix the "Employer Skills Survey" entry (A11): typo "opf" -> "of"
$ws.Range("A11").Value = "Employer Skills Survey including hard-to-fill and skill-shortage vacancies, employer skills needed in next 12 months, and percent of employers providing training"

# --- Update the saved selection to A12 (matches the final authored state)
$ws.Range("A12").Select()
